$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.756.39"
$ws.Range("E2").Value = "  +2.54%  "

$ws.Range("D3").Value = "1.707.00"
$ws.Range("E3").Value = "  +1.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.40"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9983"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3743"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.25"
$ws.Range("E8").Value = "  +3.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3445"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.195"
$ws.Range("E10").Value = "  +0.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07490"
$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9983"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.99"
$ws.Range("E13").Value = "  +3.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.253"
$ws.Range("E14").Value = "  +2.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.970"
$ws.Range("E15").Value = "  +2.88%  "

$ws.Range("D16").Value = "1.707.97"
$ws.Range("E16").Value = "  +1.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001130"
$ws.Range("E17").Value = "  +2.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06734"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9985"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "84.61"
$ws.Range("E20").Value = "  +3.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.20"
$ws.Range("E21").Value = "  +4.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.337"
$ws.Range("E22").Value = "  +3.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.04"
$ws.Range("E23").Value = "  +8.65%  "

$ws.Range("D24").Value = "24.716.57"
$ws.Range("E24").Value = "  +2.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.436"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.773"
$ws.Range("E26").Value = "  +3.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.27"
$ws.Range("E27").Value = "  +3.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.65"
$ws.Range("E28").Value = "  -1.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "131.51"
$ws.Range("E29").Value = "  +3.31%  "

$ws.Range("D30").Value = "1.895.32"
$ws.Range("E30").Value = "  +1.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.182"
$ws.Range("E31").Value = "  +20.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.806"
$ws.Range("E32").Value = "  +6.25%  "

$ws.Range("E33").Value = "  +2.75%  "

$ws.Range("E34").Value = "  +2.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08849"
$ws.Range("E35").Value = "  +4.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.71"
$ws.Range("E36").Value = "  +10.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.551"
$ws.Range("E37").Value = "  +3.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06597"
$ws.Range("E38").Value = "  +2.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.995"
$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02394"
$ws.Range("E40").Value = "  +1.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2230"
$ws.Range("E41").Value = "  +5.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.275"
$ws.Range("E42").Value = "  +0.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6461"
$ws.Range("E43").Value = "  +4.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9976"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.93"
$ws.Range("E45").Value = "  +5.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6132"
$ws.Range("E46").Value = "  +2.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.819"
$ws.Range("E47").Value = "  +0.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.124"
$ws.Range("E48").Value = "  +4.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.66"
$ws.Range("E49").Value = "  +2.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07296"
$ws.Range("E50").Value = "  +1.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.46"
